$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$win.ScrollRow = 10
$win.ScrollColumn = 1
Write-Output ("ScrollRow after unfreeze: " + $win.ScrollRow())
$win.FreezePanes = $true
Write-Output ("SplitRow after refreeze: " + $win.SplitRow())
Write-Output ("ScrollRow after refreeze: " + $win.ScrollRow())
$ws.Range("B21").Select()
Write-Output "done"
